# GHI_2024-02-23.xlsx refresh: a newer solar-irradiance forecast run
# (new lat/lon fix + new sunrise/sunset + new clear/cloudy-sky GHI/DNI/DHI
# numbers) replaces the previous run's data on both the "Daily" and
# "Hourly" sheets.

$wb = $excel.ActiveWorkbook
$wsDaily  = $wb.Worksheets.Item("Daily")
$wsHourly = $wb.Worksheets.Item("Hourly")

# --- Daily sheet: row 2 lat/lon + daily totals ---
$wsDaily.Range("A2").Value = 47.2229
$wsDaily.Range("B2").Value = 24.7244
$wsDaily.Range("G2").Value = 3426.32
$wsDaily.Range("H2").Value = 6752.65
$wsDaily.Range("I2").Value = 798
$wsDaily.Range("J2").Value = 867.67
$wsDaily.Range("L2").Value = 867.67
# (K2 stays 0 - unchanged)

# --- Hourly sheet: lat/lon are repeated on every data row (2-25) ---
$wsHourly.Range("A2:A25").Value = 47.2229
$wsHourly.Range("B2:B25").Value = 24.7244

# --- Hourly sheet: per-hour clear/cloudy-sky irradiance updates ---
# (only the daylight hours, rows 9-19, actually change - the rest stay 0)
$wsHourly.Range("H9").Value  = 22.52
$wsHourly.Range("I9").Value  = 139.29
$wsHourly.Range("J9").Value  = 18.41
$wsHourly.Range("K9").Value  = 5.97
$wsHourly.Range("M9").Value  = 5.97

$wsHourly.Range("H10").Value = 151.45
$wsHourly.Range("I10").Value = 512.33
$wsHourly.Range("J10").Value = 57.33
$wsHourly.Range("K10").Value = 39.05
$wsHourly.Range("M10").Value = 39.05

$wsHourly.Range("H11").Value = 299.63
$wsHourly.Range("I11").Value = 679.22
$wsHourly.Range("J11").Value = 78.05
$wsHourly.Range("K11").Value = 74.91
$wsHourly.Range("M11").Value = 74.91

$wsHourly.Range("H12").Value = 423.23
$wsHourly.Range("I12").Value = 763.52
$wsHourly.Range("J12").Value = 90.47
$wsHourly.Range("K12").Value = 105.81
$wsHourly.Range("M12").Value = 105.81

$wsHourly.Range("H13").Value = 505.32
$wsHourly.Range("I13").Value = 806.22
$wsHourly.Range("J13").Value = 97.40000000000001
$wsHourly.Range("K13").Value = 126.33
$wsHourly.Range("M13").Value = 126.33

$wsHourly.Range("H14").Value = 536.65
$wsHourly.Range("I14").Value = 820.6
$wsHourly.Range("J14").Value = 99.84
$wsHourly.Range("K14").Value = 134.44
$wsHourly.Range("M14").Value = 134.44

$wsHourly.Range("H15").Value = 513.92
$wsHourly.Range("I15").Value = 810.33
$wsHourly.Range("J15").Value = 98.06
$wsHourly.Range("K15").Value = 129.2
$wsHourly.Range("M15").Value = 129.2

$wsHourly.Range("H16").Value = 439.5
$wsHourly.Range("I16").Value = 772.74
$wsHourly.Range("J16").Value = 91.89
$wsHourly.Range("K16").Value = 110.27
$wsHourly.Range("M16").Value = 110.27

$wsHourly.Range("H17").Value = 321.67
$wsHourly.Range("I17").Value = 696.72
$wsHourly.Range("J17").Value = 80.47
$wsHourly.Range("K17").Value = 81.77
$wsHourly.Range("M17").Value = 81.77

$wsHourly.Range("H18").Value = 175.87
$wsHourly.Range("I18").Value = 548.72
$wsHourly.Range("J18").Value = 61.47
$wsHourly.Range("K18").Value = 48.68
$wsHourly.Range("M18").Value = 48.68

$wsHourly.Range("H19").Value = 36.56
$wsHourly.Range("I19").Value = 202.96
$wsHourly.Range("J19").Value = 24.61
$wsHourly.Range("K19").Value = 11.24
$wsHourly.Range("M19").Value = 11.24

# --- Sunrise / sunset timestamps (columns E / F), repeated on every row ---
# Every cell referencing the old shared strings is rewritten so the stale
# "2024-02-23T07:16:04" / "2024-02-23T18:02:46" text is fully retired and
# the new timestamps are shared consistently across both sheets.
$wsDaily.Range("E2").Value  = "2024-02-23T07:13:08"
$wsDaily.Range("F2").Value  = "2024-02-23T17:56:35"
$wsHourly.Range("E2:E25").Value = "2024-02-23T07:13:08"
$wsHourly.Range("F2:F25").Value = "2024-02-23T17:56:35"

# --- View state: the refreshed workbook now opens on the Daily sheet ---
$wsDaily.Activate()
